# Fruta / hortaliza, semanal
# Inserts the latest weekly Papaya price observation (Vega Modelo de Temuco,
# La Araucanía) as a new row 106, pushing the previously existing rows
# 106..127 down to 107..128.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 106 (shifts 106:127 -> 107:128)
$ws.Rows.Item(106).Insert()

# Populate the new row 106 with this week's data point
$ws.Cells.Item(106, 1).Value = 10
$ws.Cells.Item(106, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(106, 3).Value = "La Araucanía"
$ws.Cells.Item(106, 4).Value = 45211
$ws.Cells.Item(106, 5).Value = 9
$ws.Cells.Item(106, 6).Value = "Fruta"
$ws.Cells.Item(106, 7).Value = 100108
$ws.Cells.Item(106, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(106, 9).Value = 100108004
$ws.Cells.Item(106, 10).Value = "Papaya"
$ws.Cells.Item(106, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(106, 12).Value = "Primera"
$ws.Cells.Item(106, 13).Value = 305
$ws.Cells.Item(106, 14).Value = 24000
$ws.Cells.Item(106, 15).Value = 25000
$ws.Cells.Item(106, 16).Value = 24180
$ws.Cells.Item(106, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(106, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(106, 19).Value = 2418
$ws.Cells.Item(106, 20).Value = 10
